$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.941.21"
$ws.Range("E2").Value = "  +2.57%  "

$ws.Range("D3").Value = "3.230.03"
$ws.Range("E3").Value = "  +6.53%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "'577.51"
$ws.Range("E5").Value = "  +4.48%  "

$ws.Range("D6").Value = "'150.07"
$ws.Range("E6").Value = "  +8.08%  "

$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("D8").Value = "3.219.76"
$ws.Range("E8").Value = "  +6.53%  "

$ws.Range("D9").Value = "'0.510"
$ws.Range("E9").Value = "  +5.34%  "

$ws.Range("D10").Value = "'7.08"
$ws.Range("E10").Value = "  +9.59%  "

$ws.Range("E11").Value = "  +6.08%  "

$ws.Range("D12").Value = "'0.484"
$ws.Range("E12").Value = "  +5.63%  "

$ws.Range("D13").Value = "'37.83"
$ws.Range("E13").Value = "  +5.43%  "

$ws.Range("D14").Value = "'0.0000231"
$ws.Range("E14").Value = "  +6.32%  "

$ws.Range("D15").Value = "3.741.55"
$ws.Range("E15").Value = "  +6.82%  "

$ws.Range("D16").Value = "66.084.97"
$ws.Range("E16").Value = "  +2.69%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.234.28"
$ws.Range("E17").Value = "  +6.43%  "

$ws.Range("B18").Value = "BitcoinCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D18").Value = "'537.46"
$ws.Range("E18").Value = "  +11.40%  "

$ws.Range("E19").Value = "  +3.08%  "

$ws.Range("D20").Value = "'7.07"
$ws.Range("E20").Value = "  +7.00%  "

$ws.Range("D21").Value = "'14.42"
$ws.Range("E21").Value = "  +7.11%  "

$ws.Range("D22").Value = "'0.739"
$ws.Range("E22").Value = "  +8.80%  "

$ws.Range("D23").Value = "'7.76"
$ws.Range("E23").Value = "  +10.19%  "

$ws.Range("D24").Value = "'13.39"
$ws.Range("E24").Value = "  +6.85%  "

$ws.Range("D25").Value = "'80.55"
$ws.Range("E25").Value = "  +3.13%  "

$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.06%  "

$ws.Range("D27").Value = "'9.26"
$ws.Range("E27").Value = "  +20.72%  "

$ws.Range("D28").Value = "'2.95"
$ws.Range("E28").Value = "  +8.82%  "

$ws.Range("D30").Value = "'27.44"
$ws.Range("E30").Value = "  +6.93%  "

$ws.Range("B31").Value = "FirstDigitalUSD"
$ws.Range("C31").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "  +0.03%  "

$ws.Range("B32").Value = "Stacks"
$ws.Range("C32").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D32").Value = "'2.72"
$ws.Range("E32").Value = "  +5.82%  "

$ws.Range("D33").Value = "'1.17"
$ws.Range("E33").Value = "  +6.12%  "

$ws.Range("D34").Value = "'560.91"
$ws.Range("E34").Value = "  +1.23%  "

$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").Value = "'6.30"
$ws.Range("E35").Value = "  +7.11%  "

$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").Value = "'5.58"
$ws.Range("E36").Value = "  +4.50%  "

$ws.Range("D37").Value = "'54.62"
$ws.Range("E37").Value = "  +4.28%  "

$ws.Range("D38").Value = "'0.0450"
$ws.Range("E38").Value = "  +8.62%  "

$ws.Range("D39").Value = "'0.0853"
$ws.Range("E39").Value = "  +7.66%  "

$ws.Range("E40").Value = "  +6.24%  "

$ws.Range("D41").Value = "3.175.81"
$ws.Range("E41").Value = "  +10.64%  "

$ws.Range("D42").Value = "'2.89"
$ws.Range("E42").Value = "  +8.03%  "

$ws.Range("D43").Value = "'8.50"
$ws.Range("E43").Value = "  +3.92%  "

$ws.Range("D44").Value = "'0.282"
$ws.Range("E44").Value = "  +17.89%  "

$ws.Range("E45").Value = "  +12.02%  "

$ws.Range("D46").Value = "'26.22"
$ws.Range("E46").Value = "  +7.30%  "

$ws.Range("D48").Value = "0.0₃0551"
$ws.Range("E48").Value = "  +5.34%  "

$ws.Range("D49").Value = "'124.90"
$ws.Range("E49").Value = "  +4.44%  "

$ws.Range("E50").Value = "  +3.73%  "

$ws.Range("D51").Value = "'2.18"
$ws.Range("E51").Value = "  +8.56%  "
